$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update region names (text changes) - apply to ALL cells referencing the old text
for ($r = 2; $r -le 81; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -eq "Asturias, Principado de") { $cell.Value = "Asturias" }
    elseif ($val -eq "Baleares, Illes") { $cell.Value = "Baleares" }
    elseif ($val -eq "Madrid, Comunidad de") { $cell.Value = "Madrid" }
    elseif ($val -eq "Murcia, Región de") { $cell.Value = "Murcia" }
    elseif ($val -eq "Navarra, Comunidad Foral de") { $cell.Value = "Navarra" }
    elseif ($val -eq "Rioja, La") { $cell.Value = "La Rioja" }
}

# Row 81: Melilla -> ñ
$ws.Range("A81").Value = "ñ"

# Row heights: rows 2-41 change from 16.5 to 18.75
for ($r = 2; $r -le 41; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
}

# Font color: data cells (B:E, rows 2-81) theme color -> explicit black RGB
$ws.Range("B2:E81").Font.Color = 0
